$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Production Units"

# New data row (row 2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Sample Production Unit"
$ws.Range("C2").Value = "Delhi, India"
$ws.Range("D2").Value = "active"

# E2 needs to hold the text "0" (not a number) while keeping the
# currency-formatted style inherited from column E. Writing the text
# directly triggers Excel's numeric auto-detection (and a quote-prefix
# style), so instead compute it as a formula and convert it in place to a
# literal value via copy / paste-special (values only).
$e2 = $ws.Cells.Item(2, 5)
$e2.Formula = '=TEXT(0,"0")'
$e2.Copy()
$e2.PasteSpecial(-4163)

$ws.Range("F2").Value = "2025-05-05T08:47:19.803Z"

# Column width updates (the first four columns were one merged <col>
# block at width 10; columns B and C now get their own custom widths,
# which naturally splits that block into individual <col> entries).
$ws.Columns.Item(2).ColumnWidth = 23.17
$ws.Columns.Item(3).ColumnWidth = 13.17
$ws.Columns.Item(6).ColumnWidth = 25.17
